# "no more b&b error"
# Updates the results sheet: fixes row 2/3 numbers and appends the
# remaining 11 solver runs (rows 4-14) that were missing before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking piece of text (e.g. "55.119...") to be
# stored as a genuine text/shared-string cell instead of being silently
# re-interpreted as a number. We do this by writing a formula that
# evaluates to the literal text, then collapsing the formula down to its
# value with a copy / paste-special(values-only) so no formula is left
# behind and no cell formatting changes.
function Set-TextValue($rangeAddr, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell = $ws.Range($rangeAddr)
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# ---- fix existing rows 2 and 3 ----
$ws.Range("B2").Value = 0.0
$ws.Range("E2").Value = 0.0
$ws.Range("F2").Value = 0.0

$ws.Range("B3").Value = 0.103
$ws.Range("E3").Value = 136.98304626489937
$ws.Range("F3").Value = 0.00008927337734168181

# ---- append the missing solver runs (rows 4-14) ----
$rows = @(
    @{ A = "10_ulysses_6.tsp"; B = 0.235;               C = "55.11939124322688";  D = "[6, 9][2, 3][5][10][1, 4, 8][7]";                                                  E = 55.11939124322688;  F = 0.0 },
    @{ A = "10_ulysses_9.tsp"; B = 1.268;               C = "33.29189782877749";  D = "[2, 3][5][9][4][1][6][10][8][7]";                                                  E = 33.29189782877749;  F = 0.0 },
    @{ A = "14_burma_3.tsp";   B = 0.543;               C = "93.3899872599682";   D = "[1, 8, 9, 10, 11][2, 12, 13, 14][3, 4, 5, 6, 7]";                                  E = 93.3899872599682;   F = 0.0 },
    @{ A = "14_burma_6.tsp";   B = 0.8489999771118164;  C = "42.74062354260174";  D = "[2, 13][3, 4, 5][9, 10][7][1, 8, 11][6, 12, 14]";                                  E = 42.74062354260174;  F = 0.0 },
    @{ A = "14_burma_9.tsp";   B = 30.075999975204468;  C = "20.762438566071065"; D = "[3, 4][2][6, 12][1, 8][9, 11][10][5][13, 14][7]";                                  E = 0.0;                F = 0.9999999999951837 },
    @{ A = "22_ulysses_3.tsp"; B = 6.430999994277954;   C = "358.6368286225183";  D = "[1, 2, 3, 4, 8, 16, 17, 18][7, 10, 12, 13, 14, 19, 20, 22][5, 6, 9, 11, 15, 21]";   E = 358.6116071877678;  F = 0.00007032583588074738 },
    @{ A = "22_ulysses_6.tsp"; B = 30.07200002670288;   C = "145.4445609954842";  D = "[10, 12, 13, 16][5, 6, 14, 15][1, 2, 7, 17][19, 20, 21][9, 11][3, 4, 8, 18, 22]";   E = 51.615014703772715; F = 0.6451224139936056 },
    @{ A = "22_ulysses_9.tsp"; B = 30.128000020980835;  C = "94.6846493760953";   D = "[17, 21][1, 9, 20][3, 15][5, 7, 13][6, 12, 14][4, 8, 18][11][10, 19][2, 16, 22]";   E = 0.0;                F = 0.9999999999989438 },
    @{ A = "26_eil_3.tsp";     B = 30.06500005722046;   C = "2743.859786570546";  D = "[4, 13, 18, 19, 24, 25][1, 5, 6, 7, 8, 14, 17, 22, 23, 26][2, 3, 9, 10, 11, 12, 15, 16, 20, 21]"; E = 1221.3575832029364; F = 0.5548760949153587 },
    @{ A = "26_eil_6.tsp";     B = 30.09000015258789;   C = "1249.5716158597847"; D = "[2, 6, 16, 21][1, 3, 7, 20, 22, 26][11, 14, 18, 25][4, 13, 19][5, 9, 10, 12, 15][8, 17, 23, 24]";   E = 290.5216572111246;  F = 0.767502995807644 },
    @{ A = "26_eil_9.tsp";     B = 30.1010000705719;    C = "1085.7153874585001"; D = "[4, 26][2, 21][1, 3, 10, 16, 22][6, 8, 24][9, 19][5, 18, 20][7, 11, 25][12, 14, 15, 17][13, 23]";   E = 0.0;                F = 0.9999999999999079 }
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    Set-TextValue $ws.Cells.Item($r, 3).Address() $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $r = $r + 1
}
